$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh reshuffles which market-day record sits on which row:
# the data for columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg) and T (Kg / unidad) move to
# a different row; columns A,B,C,E,F,G,H,I,J,K stay put since they are the
# same for every record in this sub-workbook.
#
# destRow = srcRow-that-supplies-destRow's-new-data
$rowMap = @{
    2 = 3;  3 = 24; 4 = 11; 5 = 5;  6 = 21; 7 = 31; 8 = 26; 9 = 34; 10 = 25;
    11 = 36; 12 = 27; 13 = 12; 14 = 38; 15 = 39; 16 = 30; 17 = 8; 18 = 18;
    19 = 4; 20 = 32; 21 = 20; 22 = 16; 23 = 29; 24 = 2; 25 = 22; 26 = 33;
    27 = 19; 28 = 15; 29 = 7; 30 = 28; 31 = 10; 32 = 17; 33 = 9; 34 = 23;
    35 = 14; 36 = 6; 37 = 13; 38 = 35; 39 = 37
}

$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D, L, M, N, O, P, Q, R, S, T

$firstRow = 2
$lastRow = 39

# Snapshot every cell we might touch before writing anything, so the shuffle
# reads consistent "before" data regardless of write order.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
